$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Info" after the last existing sheet
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Info"

$ws.Range("A1").Value = "Page count"
$ws.Range("B1").Value = 2

$ws.Range("A2").Value = "Page list"

# "29"/"30" are page numbers stored as text, not numeric values -
# force the cells to Text format first so Excel doesn't coerce them.
$ws.Range("B2:C2").NumberFormat = "@"
$ws.Range("B2").Value = "29"
$ws.Range("C2").Value = "30"
